$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.849.94"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.639.38"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0793"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.45%  "
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").Value = "1.865.32"
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").Value = "1.634.10"
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.565"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").Value = "25.860.10"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.41%  "
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("E23").Value = "  +1.98%  "
$ws.Range("E24").Value = "  +3.87%  "
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.94%  "
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("E28").Value = "  +1.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.38%  "
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("E32").Value = "  +1.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.30%  "
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.911"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("D37").Value = "1.134.09"
$ws.Range("E37").Value = "  +0.61%  "
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.546"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E42").Value = "  +1.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.809"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.98%  "
$ws.Range("D45").Value = "1.774.65"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("E46").Value = "  +1.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0504"
$ws.Range("D49").Style = "Normal"
$ws.Range("E50").Value = "  +4.42%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0963"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.81%  "
